$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new text value (values are percentages / decimal numbers that
# must be written back as literal text, matching the original inlineStr cells)
$updates = @{
    'D2' = '303.83'
    'E2' = '5.47%'
    'D3' = '35.07'
    'E3' = '13.05%'
    'D4' = '5.177'
    'E4' = '5.09%'
    'D5' = '0.07789'
    'E5' = '6.15%'
    'D6' = '2.383'
    'E6' = '8.68%'
    'D7' = '8.052'
    'E7' = '4.16%'
    'D8' = '3.973'
    'D9' = '0.9267'
    'E9' = '2.65%'
    'D10' = '0.09979'
    'E10' = '9.00%'
    'D11' = '0.1831'
    'E11' = '8.81%'
    'D12' = '0.08642'
    'E12' = '4.41%'
    'D13' = '0.03325'
    'E13' = '6.57%'
    'D14' = '0.09912'
    'E14' = '-0.31%'
    'D15' = '0.001507'
    'E15' = '0.84%'
    'D16' = '0.005724'
    'E16' = '0.10%'
    'E17' = '-1.54%'
    'D18' = '2.146'
    'E18' = '4.69%'
    'D19' = '0.3377'
    'E19' = '1.34%'
    'D20' = '0.1326'
    'E20' = '2.60%'
    'D21' = '4.327'
    'E21' = '3.66%'
    'D22' = '0.2387'
    'E22' = '13.53%'
    'D23' = '0.04572'
    'E23' = '0.86%'
    'D24' = '0.001217'
    'E24' = '0.65%'
    'E25' = '7.07%'
    'D26' = '0.0001298'
    'E26' = '-0.28%'
    'D27' = '0.0003701'
    'E27' = '8.92%'
    'D39' = '0.01774'
    'E39' = '13.20%'
    'D40' = '0.04791'
    'E40' = '7.96%'
    'D41' = '0.007728'
    'E41' = '4.52%'
    'E42' = '6.55%'
    'D43' = '0.007166'
    'E43' = '-20.27%'
    'D44' = '0.002239'
    'E44' = '0.27%'
    'D45' = '0.009549'
    'E45' = '4.79%'
    'D46' = '0.00006106'
    'E46' = '-0.28%'
    'D47' = '0.00000000750'
    'E47' = '-0.07%'
    'E48' = '20.29%'
    'E49' = '-0.11%'
    'D50' = '0.00002100'
    'E50' = '-0.07%'
    'D51' = '0.0002000'
    'E51' = '-0.07%'
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    # Force text storage so Excel does not reinterpret "303.83" or "5.47%"
    # as a number/percentage, then strip the temporary format delta so the
    # cell style matches the unmodified neighbours exactly.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
    $cell.ClearFormats()
}
